# Apply targeted odds/value updates to Sheet1 as described by the commit diff.
# Each assignment below updates a single cell's numeric value in place,
# matching the before/after cell values from the canonical OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N3").Value = 1.75
$ws.Range("O3").Value = 2.05
$ws.Range("L5").Value = 1.36
$ws.Range("M5").Value = 3
$ws.Range("T7").Value = 11
$ws.Range("AH7").Value = 19
$ws.Range("J8").Value = 1.04
$ws.Range("K8").Value = 13
$ws.Range("N8").Value = 1.85
$ws.Range("O8").Value = 2
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 2.95
$ws.Range("I15").Value = 2.4
$ws.Range("L15").Value = 1.42
$ws.Range("M15").Value = 2.45
$ws.Range("N15").Value = 2.22
$ws.Range("O15").Value = 1.52
$ws.Range("P15").Value = 1.47
$ws.Range("Q15").Value = 2.32
$ws.Range("R15").Value = 1.91
$ws.Range("S15").Value = 1.7
$ws.Range("T15").Value = 7.6
$ws.Range("U15").Value = 14.5
$ws.Range("V15").Value = 11
$ws.Range("W15").Value = 37
$ws.Range("X15").Value = 29
$ws.Range("Y15").Value = 45
$ws.Range("Z15").Value = 7
$ws.Range("AA15").Value = 5.8
$ws.Range("AB15").Value = 16
$ws.Range("AC15").Value = 90
$ws.Range("AD15").Value = 900
$ws.Range("AE15").Value = 6.4
$ws.Range("AF15").Value = 10.5
$ws.Range("AG15").Value = 9.75
$ws.Range("AH15").Value = 25
$ws.Range("G16").Value = 2.27
$ws.Range("H16").Value = 2.9
$ws.Range("K16").Value = 5.9
$ws.Range("N16").Value = 2.22
$ws.Range("P16").Value = 1.45
$ws.Range("Q16").Value = 2.55
$ws.Range("R16").Value = 1.87
$ws.Range("T16").Value = 6.7
$ws.Range("U16").Value = 10.5
$ws.Range("Z16").Value = 5.9
$ws.Range("AA16").Value = 5.6
$ws.Range("AB16").Value = 14.5
$ws.Range("AI16").Value = 35
$ws.Range("G18").Value = 1.78
$ws.Range("H18").Value = 3.85
$ws.Range("I18").Value = 3.7
$ws.Range("S18").Value = 2.05
$ws.Range("T18").Value = 8.75
$ws.Range("U18").Value = 9.5
$ws.Range("V18").Value = 8.25
$ws.Range("W18").Value = 15
$ws.Range("X18").Value = 13
$ws.Range("AA18").Value = 7.7
$ws.Range("AB18").Value = 14
$ws.Range("AC18").Value = 55
$ws.Range("AE18").Value = 13.5
$ws.Range("AF18").Value = 22
$ws.Range("AH18").Value = 55
$ws.Range("AI18").Value = 30
$ws.Range("G20").Value = 3.4
$ws.Range("I20").Value = 2.2
$ws.Range("K20").Value = 8.5
$ws.Range("O20").Value = 1.67
$ws.Range("T20").Value = 9.5
$ws.Range("U20").Value = 17
$ws.Range("AJ20").Value = 29
$ws.Range("G21").Value = 2.4
$ws.Range("Y21").Value = 26
$ws.Range("Z21").Value = 11
$ws.Range("AE21").Value = 10
$ws.Range("L22").Value = 1.44
$ws.Range("M22").Value = 2.63
$ws.Range("N22").Value = 2.4
$ws.Range("O22").Value = 1.5
$ws.Range("G23").Value = 2.15
$ws.Range("H23").Value = 3
$ws.Range("I23").Value = 3.75
$ws.Range("O23").Value = 1.67
$ws.Range("W23").Value = 19
$ws.Range("AF23").Value = 19
$ws.Range("H24").Value = 3
$ws.Range("I24").Value = 5
$ws.Range("K24").Value = 8.5
$ws.Range("O24").Value = 1.67
$ws.Range("AD24").Value = 301
$ws.Range("J25").Value = 1.06
$ws.Range("K25").Value = 10
$ws.Range("L25").Value = 1.29
$ws.Range("M25").Value = 3.5
$ws.Range("N25").Value = 2
$ws.Range("O25").Value = 1.85
$ws.Range("Z25").Value = 9.5
$ws.Range("AB25").Value = 13
$ws.Range("AE25").Value = 8
$ws.Range("J26").Value = 1.05
$ws.Range("K26").Value = 11
$ws.Range("L26").Value = 1.25
$ws.Range("M26").Value = 3.75
$ws.Range("N26").Value = 1.93
$ws.Range("O26").Value = 1.93
$ws.Range("G30").Value = 1.93
$ws.Range("I30").Value = 3.6
$ws.Range("N30").Value = 1.91
$ws.Range("T30").Value = 5.8
$ws.Range("U30").Value = 7.5
$ws.Range("W30").Value = 13.5
$ws.Range("Z30").Value = 9
$ws.Range("AC30").Value = 50
$ws.Range("AE30").Value = 8.75
$ws.Range("AF30").Value = 16
$ws.Range("AG30").Value = 10.25
$ws.Range("AH30").Value = 40
$ws.Range("AI30").Value = 26
$ws.Range("H32").Value = 5.75
$ws.Range("I32").Value = 10
$ws.Range("T32").Value = 8.5
$ws.Range("AA32").Value = 12
$ws.Range("AH32").Value = 126
$ws.Range("AI32").Value = 67
$ws.Range("H33").Value = 6.25
$ws.Range("I33").Value = 15
$ws.Range("K33").Value = 15
$ws.Range("T33").Value = 8.5
$ws.Range("Y33").Value = 34
$ws.Range("AJ33").Value = 67
$ws.Range("G35").Value = 1.75
$ws.Range("H35").Value = 3.7
$ws.Range("J35").Value = 1.04
$ws.Range("K35").Value = 9
$ws.Range("L35").Value = 1.22
$ws.Range("M35").Value = 4
$ws.Range("N35").Value = 1.75
$ws.Range("O35").Value = 2.05
$ws.Range("P35").Value = 1.33
$ws.Range("Q35").Value = 3.25
$ws.Range("R35").Value = 1.73
$ws.Range("S35").Value = 2
$ws.Range("T35").Value = 8
$ws.Range("U35").Value = 9
$ws.Range("V35").Value = 8.5
$ws.Range("Y35").Value = 23
$ws.Range("Z35").Value = 12
$ws.Range("AA35").Value = 7.5
$ws.Range("AC35").Value = 41
$ws.Range("AD35").Value = 500
$ws.Range("AE35").Value = 13
$ws.Range("AF35").Value = 23
$ws.Range("G36").Value = 2.6
$ws.Range("I36").Value = 2.38
$ws.Range("J36").Value = 21
$ws.Range("K36").Value = 1.03
$ws.Range("AC36").Value = 26
$ws.Range("AE36").Value = 15
$ws.Range("AF36").Value = 17
$ws.Range("AI36").Value = 17
$ws.Range("G37").Value = 1.35
$ws.Range("I37").Value = 8.5
$ws.Range("J37").Value = 1.05
$ws.Range("K37").Value = 8
$ws.Range("L37").Value = 1.23
$ws.Range("M37").Value = 3.75
$ws.Range("N37").Value = 1.7
$ws.Range("O37").Value = 2.05
$ws.Range("P37").Value = 1.35
$ws.Range("Q37").Value = 2.92
$ws.Range("R37").Value = 2
$ws.Range("S37").Value = 1.72
$ws.Range("T37").Value = 6.8
$ws.Range("W37").Value = 8.5
$ws.Range("X37").Value = 11.25
$ws.Range("Y37").Value = 28
$ws.Range("Z37").Value = 8
$ws.Range("AA37").Value = 9
$ws.Range("AB37").Value = 21
$ws.Range("AC37").Value = 110
$ws.Range("AD37").Value = 900
$ws.Range("AE37").Value = 20
$ws.Range("AF37").Value = 55
$ws.Range("AG37").Value = 26
$ws.Range("AJ37").Value = 90
$ws.Range("G38").Value = 1.38
$ws.Range("I38").Value = 6.7
$ws.Range("P38").Value = 1.33
$ws.Range("Q38").Value = 3.05
$ws.Range("R38").Value = 1.93
$ws.Range("S38").Value = 1.78
$ws.Range("T38").Value = 7.2
$ws.Range("U38").Value = 6.6
$ws.Range("Y38").Value = 27
$ws.Range("AB38").Value = 20
$ws.Range("AC38").Value = 90
$ws.Range("AD38").Value = 700
$ws.Range("AI38").Value = 75

Write-Output "Applied 199 cell updates."
